# Excel COM-interop script to apply the diff between before.xlsx and after.xlsx
# Commit message: Regeneration of en/fr/es country data books after child program paras updated

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 27 "Programmes pour les enfants": update base assumption cells (rows 2-53).
# All downstream cells (rows 56-163, the x0.9 / x1.05 scenarios) are driven by
# formulas referencing these base cells, so they recalc automatically.
# ---------------------------------------------------------------------------
$ws27 = $wb.Worksheets.Item("Programmes pour les enfants")

$ws27.Range("F2").Value = 0.39473684210526322
$ws27.Range("G2").Value = 0.39473684210526322
$ws27.Range("H2").Value = 0.39473684210526322
$ws27.Range("F3").Value = 0.30769230769230765
$ws27.Range("G3").Value = 0.30769230769230765
$ws27.Range("H3").Value = 0.30769230769230765
$ws27.Range("F18").Value = 0.7
$ws27.Range("F20").Value = 0.84
$ws27.Range("D21").Value = 0.28260869565217389
$ws27.Range("F21").Value = 0
$ws27.Range("F22").Value = 0
$ws27.Range("D23").Value = 0.28260869565217389
$ws27.Range("F23").Value = 0
$ws27.Range("F24").Value = 0
$ws27.Range("D25").Value = 0.28260869565217389
$ws27.Range("F25").Value = 0
$ws27.Range("F26").Value = 0
$ws27.Range("F27").Value = 1
$ws27.Range("F28").Value = 0
$ws27.Range("F29").Value = 0
$ws27.Range("F30").Value = 1
$ws27.Range("F31").Value = 0
$ws27.Range("F32").Value = 0
$ws27.Range("F33").Value = 1
$ws27.Range("F34").Value = 0
$ws27.Range("F35").Value = 0
$ws27.Range("F36").Value = 1
$ws27.Range("F37").Value = 0
$ws27.Range("F38").Value = 0
$ws27.Range("F39").Value = 1
$ws27.Range("F40").Value = 0
$ws27.Range("F41").Value = 0
$ws27.Range("F42").Value = 0.3
$ws27.Range("F43").Value = 0.5
$ws27.Range("F44").Value = 0.65
$ws27.Range("F45").Value = 0.3
$ws27.Range("F46").Value = 0.49
$ws27.Range("F47").Value = 0.52
$ws27.Range("F48").Value = 0.88
$ws27.Range("D49").Value = 0.78409090909090906
$ws27.Range("E49").Value = 0.78409090909090906
$ws27.Range("F49").Value = 0.78409090909090906
$ws27.Range("G49").Value = 0.78409090909090906
$ws27.Range("H49").Value = 0.78409090909090906
$ws27.Range("D50").Value = 0.88372093023255816
$ws27.Range("E50").Value = 0.88372093023255816
$ws27.Range("F50").Value = 0.88372093023255816
$ws27.Range("G50").Value = 0.88372093023255816
$ws27.Range("H50").Value = 0.88372093023255816
$ws27.Range("F51").Value = 0.86
$ws27.Range("F52").Value = 0
$ws27.Range("F53").Value = 0

# Update the view state on sheet 27: new top-left cell / selection
$ws27.Activate()
$ws27.Range("D2:H53").Select()

# ---------------------------------------------------------------------------
# Move the active/selected tab from sheet 1 to sheet 11 ("Dépendances du
# programme"), matching the new workbookView activeTab/sheetView tabSelected.
# ---------------------------------------------------------------------------
$ws11 = $wb.Worksheets.Item("Dépendances du programme")
$ws11.Activate()

